$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Folder Inventory": a new folder entry ("Enterprise-Class Networking
# in Azure") was refreshed/touched most recently, so it now sits at the top
# of the (date-descending) list. This pushes the existing rows 2-8 down by
# one row, and the old duplicate row that used to live at row 8 (the same
# "Enterprise-Class Networking in Azure" folder, with its previous, older
# timestamp) is dropped since it has been replaced by the new top row.
# ---------------------------------------------------------------------------
$wsInventory = $wb.Worksheets.Item("Folder Inventory")

# Shift existing rows 2..7 down into rows 3..8 (working from the bottom up
# so we don't clobber a row before it has been copied).
for ($r = 7; $r -ge 2; $r--) {
    $dest = $r + 1
    $wsInventory.Cells.Item($dest, 1).Value = $wsInventory.Cells.Item($r, 1).Value2
    $wsInventory.Cells.Item($dest, 2).Value = $wsInventory.Cells.Item($r, 2).Value2
    $wsInventory.Cells.Item($dest, 3).Value = $wsInventory.Cells.Item($r, 3).Value2
    $wsInventory.Cells.Item($dest, 4).Value = $wsInventory.Cells.Item($r, 4).Value2
    $wsInventory.Cells.Item($dest, 5).Value = $wsInventory.Cells.Item($r, 5).Value2
}

# Write the new top row (row 2).
$wsInventory.Cells.Item(2, 1).Value = "Enterprise-Class Networking in Azure"
$wsInventory.Cells.Item(2, 2).Value = "Enterprise-Class Networking in Azure"
$wsInventory.Cells.Item(2, 3).Value = "2025-06-12 12:35:48 +0530"
$wsInventory.Cells.Item(2, 4).Value = 1
$wsInventory.Cells.Item(2, 5).Value = "Root"

# ---------------------------------------------------------------------------
# Sheet "Metadata": regenerate timestamp + bump workflow run counter.
# ---------------------------------------------------------------------------
$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B3").Value = "2025-06-12 07:06:08 UTC"
$wsMetadata.Range("B5").Value = "'8"

# ---------------------------------------------------------------------------
# Sheet "Summary": the most recent update timestamp now matches the new
# top entry in the Folder Inventory sheet.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = "2025-06-12 12:35:48 +0530"
